# Apply the RLIe_scores.xlsx update: append " ecosystems" to the GET realm/
# biome labels in column A, and refresh the bootstrap confidence-interval
# bounds (columns F "lower" / G "upper") with the re-run values used for the
# new "GET per realm" barplot (commit: "added barplot - GET per realm").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('A2').Value = 'MT2: Supralittoral coastal systems (n = 10 ecosystems)'
$ws.Range('F2').Value = 0.54

$ws.Range('A3').Value = 'MT2: Supralittoral coastal systems (n = 10 ecosystems)'
$ws.Range('G3').Value = 0.8

$ws.Range('A4').Value = 'MT2: Supralittoral coastal systems (n = 10 ecosystems)'

$ws.Range('A5').Value = 'MT2: Supralittoral coastal systems (n = 10 ecosystems)'

$ws.Range('A6').Value = 'T1: Tropical-subtropical forests (n = 44 ecosystems)'

$ws.Range('A7').Value = 'T1: Tropical-subtropical forests (n = 44 ecosystems)'
$ws.Range('G7').Value = 0.9636363636363636

$ws.Range('A8').Value = 'T1: Tropical-subtropical forests (n = 44 ecosystems)'

$ws.Range('A9').Value = 'T1: Tropical-subtropical forests (n = 44 ecosystems)'

$ws.Range('A10').Value = 'T2: Temperate-boreal forests & woodlands (n = 4 ecosystems)'

$ws.Range('A11').Value = 'T2: Temperate-boreal forests & woodlands (n = 4 ecosystems)'

$ws.Range('A12').Value = 'T2: Temperate-boreal forests & woodlands (n = 4 ecosystems)'

$ws.Range('A13').Value = 'T2: Temperate-boreal forests & woodlands (n = 4 ecosystems)'

$ws.Range('A14').Value = 'T3: Shrublands & shrubby woodlands (n = 133 ecosystems)'
$ws.Range('F14').Value = 0.5879323308270676
$ws.Range('G14').Value = 0.7067669172932332

$ws.Range('A15').Value = 'T3: Shrublands & shrubby woodlands (n = 133 ecosystems)'
$ws.Range('F15').Value = 0.5714285714285714
$ws.Range('G15').Value = 0.6947368421052631

$ws.Range('A16').Value = 'T3: Shrublands & shrubby woodlands (n = 133 ecosystems)'
$ws.Range('F16').Value = 0.5698872180451128
$ws.Range('G16').Value = 0.6947368421052631

$ws.Range('A17').Value = 'T3: Shrublands & shrubby woodlands (n = 133 ecosystems)'
$ws.Range('F17').Value = 0.5684210526315789
$ws.Range('G17').Value = 0.6947368421052631

$ws.Range('A18').Value = 'T4: Savannas and grasslands (n = 180 ecosystems)'
$ws.Range('F18').Value = 0.8333333333333334

$ws.Range('A19').Value = 'T4: Savannas and grasslands (n = 180 ecosystems)'
$ws.Range('F19').Value = 0.8055555555555556

$ws.Range('A20').Value = 'T4: Savannas and grasslands (n = 180 ecosystems)'

$ws.Range('A21').Value = 'T4: Savannas and grasslands (n = 180 ecosystems)'
$ws.Range('G21').Value = 0.8722499999999996

$ws.Range('A22').Value = 'T5: Deserts and semi-deserts (n = 92 ecosystems)'
$ws.Range('F22').Value = 0.8913043478260869

$ws.Range('A23').Value = 'T5: Deserts and semi-deserts (n = 92 ecosystems)'
$ws.Range('F23').Value = 0.8695652173913043
$ws.Range('G23').Value = 0.9630434782608696

$ws.Range('A24').Value = 'T5: Deserts and semi-deserts (n = 92 ecosystems)'

$ws.Range('A25').Value = 'T5: Deserts and semi-deserts (n = 92 ecosystems)'
$ws.Range('F25').Value = 0.8673913043478261
$ws.Range('G25').Value = 0.9608695652173913

$ws.Range('A26').Value = 'MT2.1: Coastal shrublands and grasslands (n = 10 ecosystems)'

$ws.Range('A27').Value = 'MT2.1: Coastal shrublands and grasslands (n = 10 ecosystems)'

$ws.Range('A28').Value = 'MT2.1: Coastal shrublands and grasslands (n = 10 ecosystems)'

$ws.Range('A29').Value = 'MT2.1: Coastal shrublands and grasslands (n = 10 ecosystems)'

$ws.Range('A30').Value = 'T1.2: Tropical-subtropical dry forests and thickets (n = 44 ecosystems)'

$ws.Range('A31').Value = 'T1.2: Tropical-subtropical dry forests and thickets (n = 44 ecosystems)'

$ws.Range('A32').Value = 'T1.2: Tropical-subtropical dry forests and thickets (n = 44 ecosystems)'

$ws.Range('A33').Value = 'T1.2: Tropical-subtropical dry forests and thickets (n = 44 ecosystems)'

$ws.Range('A34').Value = 'T2.4: Warm temperate laurophyll forests (n = 4 ecosystems)'

$ws.Range('A35').Value = 'T2.4: Warm temperate laurophyll forests (n = 4 ecosystems)'

$ws.Range('A36').Value = 'T2.4: Warm temperate laurophyll forests (n = 4 ecosystems)'

$ws.Range('A37').Value = 'T2.4: Warm temperate laurophyll forests (n = 4 ecosystems)'

$ws.Range('A38').Value = 'T3.1: Seasonally dry tropical shrublands (n = 1 ecosystems)'

$ws.Range('A39').Value = 'T3.1: Seasonally dry tropical shrublands (n = 1 ecosystems)'

$ws.Range('A40').Value = 'T3.1: Seasonally dry tropical shrublands (n = 1 ecosystems)'

$ws.Range('A41').Value = 'T3.1: Seasonally dry tropical shrublands (n = 1 ecosystems)'

$ws.Range('A42').Value = 'T3.2: Seasonally dry temperate heaths and shrublands (n = 132 ecosystems)'
$ws.Range('G42').Value = 0.7060984848484843

$ws.Range('A43').Value = 'T3.2: Seasonally dry temperate heaths and shrublands (n = 132 ecosystems)'
$ws.Range('G43').Value = 0.6924621212121206

$ws.Range('A44').Value = 'T3.2: Seasonally dry temperate heaths and shrublands (n = 132 ecosystems)'

$ws.Range('A45').Value = 'T3.2: Seasonally dry temperate heaths and shrublands (n = 132 ecosystems)'

$ws.Range('A46').Value = 'T4.1: Trophic savannas (n = 38 ecosystems)'
$ws.Range('F46').Value = 0.8789473684210526

$ws.Range('A47').Value = 'T4.1: Trophic savannas (n = 38 ecosystems)'

$ws.Range('A48').Value = 'T4.1: Trophic savannas (n = 38 ecosystems)'

$ws.Range('A49').Value = 'T4.1: Trophic savannas (n = 38 ecosystems)'
$ws.Range('F49').Value = 0.8578947368421053

$ws.Range('A50').Value = 'T4.2: Pyric tussock savannas (n = 67 ecosystems)'

$ws.Range('A51').Value = 'T4.2: Pyric tussock savannas (n = 67 ecosystems)'
$ws.Range('F51').Value = 0.7432835820895523
$ws.Range('G51').Value = 0.8835820895522388

$ws.Range('A52').Value = 'T4.2: Pyric tussock savannas (n = 67 ecosystems)'

$ws.Range('A53').Value = 'T4.2: Pyric tussock savannas (n = 67 ecosystems)'
$ws.Range('F53').Value = 0.7313432835820896
$ws.Range('G53').Value = 0.8746268656716418

$ws.Range('A54').Value = 'T4.5: Temperate subhumid grasslands (n = 75 ecosystems)'
$ws.Range('G54').Value = 0.904

$ws.Range('A55').Value = 'T4.5: Temperate subhumid grasslands (n = 75 ecosystems)'
$ws.Range('F55').Value = 0.7626666666666666
$ws.Range('G55').Value = 0.885399999999999

$ws.Range('A56').Value = 'T4.5: Temperate subhumid grasslands (n = 75 ecosystems)'
$ws.Range('F56').Value = 0.7599333333333335
$ws.Range('G56').Value = 0.8826666666666667

$ws.Range('A57').Value = 'T4.5: Temperate subhumid grasslands (n = 75 ecosystems)'
$ws.Range('F57').Value = 0.7466666666666666

$ws.Range('A58').Value = 'T5.1: Semi-desert steppes (n = 15 ecosystems)'

$ws.Range('A59').Value = 'T5.1: Semi-desert steppes (n = 15 ecosystems)'

$ws.Range('A60').Value = 'T5.1: Semi-desert steppes (n = 15 ecosystems)'

$ws.Range('A61').Value = 'T5.1: Semi-desert steppes (n = 15 ecosystems)'

$ws.Range('A62').Value = 'T5.2: Thorny deserts and semi-deserts (n = 63 ecosystems)'
$ws.Range('F62').Value = 0.8603174603174604
$ws.Range('G62').Value = 0.9714285714285714

$ws.Range('A63').Value = 'T5.2: Thorny deserts and semi-deserts (n = 63 ecosystems)'
$ws.Range('G63').Value = 0.9619047619047619

$ws.Range('A64').Value = 'T5.2: Thorny deserts and semi-deserts (n = 63 ecosystems)'

$ws.Range('A65').Value = 'T5.2: Thorny deserts and semi-deserts (n = 63 ecosystems)'

$ws.Range('A66').Value = 'T5.5: Hyper-arid deserts (n = 14 ecosystems)'

$ws.Range('A67').Value = 'T5.5: Hyper-arid deserts (n = 14 ecosystems)'

$ws.Range('A68').Value = 'T5.5: Hyper-arid deserts (n = 14 ecosystems)'

$ws.Range('A69').Value = 'T5.5: Hyper-arid deserts (n = 14 ecosystems)'

$ws.Range('F70').Value = 0.7939524838012959
$ws.Range('G70').Value = 0.8483801295896328

$ws.Range('F71').Value = 0.7688984881209503

$ws.Range('F72').Value = 0.7676025917926566
$ws.Range('G72').Value = 0.82463282937365

$ws.Range('F73').Value = 0.7658747300215982
$ws.Range('G73').Value = 0.8237580993520518
